# Atualização de bases das ligas, do dia: 17-05-2024 às 13:59
#
# Row 11 and Row 12 of the "Germany Landesliga" sheet actually describe the
# same two matches but with their full records (id, teams, score, result,
# odds) swapped between the two rows. The row-index column (A) and the
# Div/Date columns (C/D) are identical between the two rows, so only B and
# F..AB (plus E, the HomeTeam) need to be written with the swapped values.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Germany Landesliga")

# --- snapshot the current ("before") values of row 11 and row 12 ---------
$b11  = $ws.Range("B11").Value2
$e11  = $ws.Range("E11").Value2
$f11  = $ws.Range("F11").Value2
$g11  = $ws.Range("G11").Value2
$h11  = $ws.Range("H11").Value2
$i11  = $ws.Range("I11").Value2
$j11  = $ws.Range("J11").Value2
$k11  = $ws.Range("K11").Value2
$l11  = $ws.Range("L11").Value2
$m11  = $ws.Range("M11").Value2
$n11  = $ws.Range("N11").Value2
$o11  = $ws.Range("O11").Value2
$p11  = $ws.Range("P11").Value2
$q11  = $ws.Range("Q11").Value2
$r11  = $ws.Range("R11").Value2
$s11  = $ws.Range("S11").Value2
$t11  = $ws.Range("T11").Value2
$u11  = $ws.Range("U11").Value2
$v11  = $ws.Range("V11").Value2
$w11  = $ws.Range("W11").Value2
$x11  = $ws.Range("X11").Value2
$y11  = $ws.Range("Y11").Value2
$z11  = $ws.Range("Z11").Value2
$aa11 = $ws.Range("AA11").Value2
$ab11 = $ws.Range("AB11").Value2

$b12  = $ws.Range("B12").Value2
$e12  = $ws.Range("E12").Value2
$f12  = $ws.Range("F12").Value2
$g12  = $ws.Range("G12").Value2
$h12  = $ws.Range("H12").Value2
$i12  = $ws.Range("I12").Value2
$j12  = $ws.Range("J12").Value2
$k12  = $ws.Range("K12").Value2
$l12  = $ws.Range("L12").Value2
$m12  = $ws.Range("M12").Value2
$n12  = $ws.Range("N12").Value2
$o12  = $ws.Range("O12").Value2
$p12  = $ws.Range("P12").Value2
$q12  = $ws.Range("Q12").Value2
$r12  = $ws.Range("R12").Value2
$s12  = $ws.Range("S12").Value2
$t12  = $ws.Range("T12").Value2
$u12  = $ws.Range("U12").Value2
$v12  = $ws.Range("V12").Value2
$w12  = $ws.Range("W12").Value2
$x12  = $ws.Range("X12").Value2
$y12  = $ws.Range("Y12").Value2
$z12  = $ws.Range("Z12").Value2
$aa12 = $ws.Range("AA12").Value2
$ab12 = $ws.Range("AB12").Value2

# --- write row 11's cells with row 12's old values ------------------------
$ws.Range("B11").Value2  = $b12
$ws.Range("E11").Value2  = $e12
$ws.Range("F11").Value2  = $f12
$ws.Range("G11").Value2  = $g12
$ws.Range("H11").Value2  = $h12
$ws.Range("I11").Value2  = $i12
$ws.Range("J11").Value2  = $j12
$ws.Range("K11").Value2  = $k12
$ws.Range("L11").Value2  = $l12
$ws.Range("M11").Value2  = $m12
$ws.Range("N11").Value2  = $n12
$ws.Range("O11").Value2  = $o12
$ws.Range("P11").Value2  = $p12
$ws.Range("Q11").Value2  = $q12
$ws.Range("R11").Value2  = $r12
$ws.Range("S11").Value2  = $s12
$ws.Range("T11").Value2  = $t12
$ws.Range("U11").Value2  = $u12
$ws.Range("V11").Value2  = $v12
$ws.Range("W11").Value2  = $w12
$ws.Range("X11").Value2  = $x12
$ws.Range("Y11").Value2  = $y12
$ws.Range("Z11").Value2  = $z12
$ws.Range("AA11").Value2 = $aa12
$ws.Range("AB11").Value2 = $ab12

# --- write row 12's cells with row 11's old values ------------------------
$ws.Range("B12").Value2  = $b11
$ws.Range("E12").Value2  = $e11
$ws.Range("F12").Value2  = $f11
$ws.Range("G12").Value2  = $g11
$ws.Range("H12").Value2  = $h11
$ws.Range("I12").Value2  = $i11
$ws.Range("J12").Value2  = $j11
$ws.Range("K12").Value2  = $k11
$ws.Range("L12").Value2  = $l11
$ws.Range("M12").Value2  = $m11
$ws.Range("N12").Value2  = $n11
$ws.Range("O12").Value2  = $o11
$ws.Range("P12").Value2  = $p11
$ws.Range("Q12").Value2  = $q11
$ws.Range("R12").Value2  = $r11
$ws.Range("S12").Value2  = $s11
$ws.Range("T12").Value2  = $t11
$ws.Range("U12").Value2  = $u11
$ws.Range("V12").Value2  = $v11
$ws.Range("W12").Value2  = $w11
$ws.Range("X12").Value2  = $x11
$ws.Range("Y12").Value2  = $y11
$ws.Range("Z12").Value2  = $z11
$ws.Range("AA12").Value2 = $aa11
$ws.Range("AB12").Value2 = $ab11
